{"js": "// Replace the multiplication-expression text runs in the practice table.\n// Each entry is [oldText, newText]; every oldText occurs exactly once in\n// the document, inside its own <w:t> run cell.\nconst replacements = [\n  [\"811\u00d74=3244\", \"814\u00d76=4884\"],\n  [\"870\u00d79=7830\", \"118\u00d72=236\"],\n  [\"153\u00d74=612\", \"843\u00d72=1686\"],\n  [\"585\u00d77=4095\", \"809\u00d76=4854\"],\n  [\"481\u00d79=4329\", \"819\u00d74=3276\"],\n  [\"711\u00d75=3555\", \"367\u00d76=2202\"],\n  [\"169\u00d76=1014\", \"358\u00d76=2148\"],\n  [\"330\u00d73=990\", \"448\u00d72=896\"],\n  [\"305\u00d74=1220\", \"600\u00d74=2400\"],\n  [\"281\u00d72=562\", \"104\u00d77=728\"],\n  [\"204\u00d75=1020\", \"710\u00d78=5680\"],\n  [\"567\u00d73=1701\", \"950\u00d72=1900\"],\n  [\"400\u00d76=2400\", \"692\u00d79=6228\"],\n  [\"666\u00d72=1332\", \"719\u00d78=5752\"],\n  [\"382\u00d72=764\", \"432\u00d78=3456\"],\n  [\"605\u00d77=4235\", \"738\u00d78=5904\"],\n  [\"105\u00d72=210\", \"390\u00d77=2730\"],\n  [\"598\u00d76=3588\", \"378\u00d74=1512\"],\n  [\"300\u00d73=900\", \"896\u00d79=8064\"],\n  [\"834\u00d75=4170\", \"559\u00d73=1677\"],\n  [\"745\u00d72=1490\", \"762\u00d74=3048\"],\n  [\"611\u00d73=1833\", \"669\u00d77=4683\"],\n  [\"637\u00d74=2548\", \"162\u00d78=1296\"],\n  [\"564\u00d74=2256\", \"548\u00d73=1644\"],\n  [\"965\u00d75=4825\", \"971\u00d73=2913\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-expression text runs in the practice table.\n# Each pair is [oldText, newText]; every oldText occurs exactly once in\n# the document.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"811\u00d74=3244\", \"814\u00d76=4884\"),\n    @(\"870\u00d79=7830\", \"118\u00d72=236\"),\n    @(\"153\u00d74=612\", \"843\u00d72=1686\"),\n    @(\"585\u00d77=4095\", \"809\u00d76=4854\"),\n    @(\"481\u00d79=4329\", \"819\u00d74=3276\"),\n    @(\"711\u00d75=3555\", \"367\u00d76=2202\"),\n    @(\"169\u00d76=1014\", \"358\u00d76=2148\"),\n    @(\"330\u00d73=990\", \"448\u00d72=896\"),\n    @(\"305\u00d74=1220\", \"600\u00d74=2400\"),\n    @(\"281\u00d72=562\", \"104\u00d77=728\"),\n    @(\"204\u00d75=1020\", \"710\u00d78=5680\"),\n    @(\"567\u00d73=1701\", \"950\u00d72=1900\"),\n    @(\"400\u00d76=2400\", \"692\u00d79=6228\"),\n    @(\"666\u00d72=1332\", \"719\u00d78=5752\"),\n    @(\"382\u00d72=764\", \"432\u00d78=3456\"),\n    @(\"605\u00d77=4235\", \"738\u00d78=5904\"),\n    @(\"105\u00d72=210\", \"390\u00d77=2730\"),\n    @(\"598\u00d76=3588\", \"378\u00d74=1512\"),\n    @(\"300\u00d73=900\", \"896\u00d79=8064\"),\n    @(\"834\u00d75=4170\", \"559\u00d73=1677\"),\n    @(\"745\u00d72=1490\", \"762\u00d74=3048\"),\n    @(\"611\u00d73=1833\", \"669\u00d77=4683\"),\n    @(\"637\u00d74=2548\", \"162\u00d78=1296\"),\n    @(\"564\u00d74=2256\", \"548\u00d73=1644\"),\n    @(\"965\u00d75=4825\", \"971\u00d73=2913\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Replacement.Text = $newText\n    $rng.Find.Forward = $true\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute(\n        $oldText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n}\n"}
